# Rename the header row (row 1) on Sheet1 to the new column names.
# "prompt title " -> "title"
# "prompt "       -> "prompt_text"
# "output format  " -> "output_format  "   (text unchanged, but re-imported)
# "department"    -> "category"
# "task type"     -> "task_type"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "task_type"
$ws.Range("C1").Value = "output_format  "
$ws.Range("D1").Value = "category"
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "prompt_text"

# Reset the view: scroll back to the top-left and select B1 (matches saved sheetView state).
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
